# BIS-1002: Fixed XLS export tests
# Add a new "Internal Assignment" column (O) to the sample type export sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for the new column, styled like the other header cells in row 4
# (bold Calibri, black).
$ws.Range("O4").Value = "Internal Assignment"
$ws.Range("O4").Font.Name = "Calibri"
$ws.Range("O4").Font.Size = 12
$ws.Range("O4").Font.Bold = $true
$ws.Range("O4").Font.Color = 0

# Data cells for the new column - plain (non-bold) text "FALSE" values for
# each of the properties listed in rows 5-9. The leading apostrophe forces
# the value to be stored as text rather than a boolean.
$ws.Range("O5").Value = "'FALSE"
$ws.Range("O6").Value = "'FALSE"
$ws.Range("O7").Value = "'FALSE"
$ws.Range("O8").Value = "'FALSE"
$ws.Range("O9").Value = "'FALSE"

$ws.Range("O5:O9").Font.Name = "Calibri"
$ws.Range("O5:O9").Font.Size = 11

# Match the updated selection left behind by the edit (O7:O9 active).
$null = $ws.Range("O7:O9").Select()
